# Generate Report for Handback
# Row 7 ("86606e37-2de0-426d-97a0-13cf99914b60.md") on both the zh-cn and
# de-de localization-status sheets finished handback processing, so its
# "Latest Target File" / "Latest Handback File" / "Latest Handback
# DateTime" / "Error Detail" columns (I, J, K, P) get populated, the same
# way rows 2-5 already are for earlier files.

$wb = $excel.ActiveWorkbook

# ---- zh-cn sheet --------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("J7").Value = "0001-01-01 00:00:00"
$wsZh.Range("K7").Value = "86606e37-2de0-426d-97a0-13cf99914b60.330b3e27d74dd55e86c1b07c8cc536b283863cd0.zh-cn.xlf"
$wsZh.Range("P7").Value = "2016-08-26 04:52:47"

$wsZh.Hyperlinks.Add(
    $wsZh.Range("I7"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/2ef24cd710d0478c4fa868447ee9b72a44acd4b4/e2e/86606e37-2de0-426d-97a0-13cf99914b60.md",
    "",
    "",
    "86606e37-2de0-426d-97a0-13cf99914b60.md"
) | Out-Null

# ---- de-de sheet ---------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("J7").Value = "86606e37-2de0-426d-97a0-13cf99914b60.330b3e27d74dd55e86c1b07c8cc536b283863cd0.de-de.xlf"
$wsDe.Range("K7").Value = "2016-08-26 04:53:16"
$wsDe.Range("P7").Value = "2016-08-26 04:52:47"

$wsDe.Hyperlinks.Add(
    $wsDe.Range("I7"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/2ef24cd710d0478c4fa868447ee9b72a44acd4b4/e2e/86606e37-2de0-426d-97a0-13cf99914b60.md",
    "",
    "",
    "86606e37-2de0-426d-97a0-13cf99914b60.md"
) | Out-Null
